# Trade #59 closed at 2026-02-17 15:43:50 - unknown UNKNOWN +0.000%
#
# This script updates the live trading results workbook to record the
# closing of trade #59:
#   - Summary sheet: roll up totals (capital, P&L, trade counts, win rate)
#   - Strategy Status sheet: roll up MarketMaking strategy stats
#   - All Trades / MarketMaking sheets: append the new trade row (#59 / row 60)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.42   # Current Capital
$summary.Range("B4").Value = 0.42      # Total P&L $
$summary.Range("B5").Value = 0.14      # Total P&L %
$summary.Range("B6").Value = 59        # Total Trades
$summary.Range("B7").Value = 19        # Winning Trades
$summary.Range("B9").Value = 32.2      # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.42     # Capital
$status.Range("D4").Value = 59         # Trades
$status.Range("E4").Value = 0.42       # P&L $
$status.Range("F4").Value = 0.42       # P&L %
$status.Range("G4").Value = 32.2       # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the new trade row to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
function Add-Trade60($ws) {
    $ws.Range("A60").Value = 59
    # Date / Time columns look like dates to type inference, so force the
    # cells to Text format first so the literal strings are preserved.
    $ws.Range("B60").NumberFormat = "@"
    $ws.Range("B60").Value = "2026-02-17"
    $ws.Range("C60").NumberFormat = "@"
    $ws.Range("C60").Value = "15:43:44"
    $ws.Range("D60").Value = "MarketMaking"
    $ws.Range("E60").Value = "DOWN"
    $ws.Range("F60").Value = 0.83
    $ws.Range("G60").Value = 0.87
    $ws.Range("H60").Value = "CLOSED"
    $ws.Range("I60").Value = 4.8193
    $ws.Range("J60").Value = 0.04
    $ws.Range("K60").Value = 100.42
    $ws.Range("L60").Value = 0
    $ws.Range("M60").Value = 0
    $ws.Range("N60").Value = 0.6
    $ws.Range("O60").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P60").Value = "early_exit"
    $ws.Range("Q60").Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade60 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade60 $marketMaking
